# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (preserve original "string" cell type)
# even when the value looks numeric, matching the source data which
# stores prices/percentages as formatted strings, not numbers.
function Set-TextValue {
    param($RangeAddr, $Text)
    $rng = $ws.Range($RangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "25.774.47"
Set-TextValue "E2" "  -2.36%  "
Set-TextValue "D3" "1.752.06"
Set-TextValue "E3" "  -4.41%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "236.98"
Set-TextValue "E5" "  -5.88%  "
Set-TextValue "E6" "  -0.16%  "
Set-TextValue "D7" "0.5070"
Set-TextValue "E7" "  -3.40%  "
Set-TextValue "D8" "41.54"
Set-TextValue "E8" "  -6.30%  "
Set-TextValue "D9" "0.2650"
Set-TextValue "E9" "  -4.63%  "
Set-TextValue "D10" "0.06161"
Set-TextValue "E10" "  -9.85%  "
Set-TextValue "D11" "1.756.24"
Set-TextValue "E11" "  -5.29%  "
Set-TextValue "D12" "15.64"
Set-TextValue "E12" "  -5.31%  "
Set-TextValue "D13" "0.06908"
Set-TextValue "E13" "  -2.63%  "
Set-TextValue "D14" "0.6029"
Set-TextValue "E14" "  -12.04%  "
Set-TextValue "D15" "4.497"
Set-TextValue "E15" "  -7.22%  "
Set-TextValue "D16" "77.23"
Set-TextValue "E16" "  -10.15%  "
Set-TextValue "E17" "  -0.04%  "
Set-TextValue "E18" "  -0.17%  "
Set-TextValue "D19" "25.781.54"
Set-TextValue "E19" "  -2.38%  "
Set-TextValue "D20" "0.000006847"
Set-TextValue "E20" "  -6.49%  "
Set-TextValue "D21" "11.70"
Set-TextValue "E21" "  -11.07%  "
Set-TextValue "D22" "1.978.05"
Set-TextValue "E22" "  -5.13%  "
Set-TextValue "D23" "4.087"
Set-TextValue "E23" "  -9.03%  "
Set-TextValue "D24" "8.247"
Set-TextValue "E24" "  -7.80%  "
Set-TextValue "D25" "5.205"
Set-TextValue "E25" "  -10.25%  "
Set-TextValue "D26" "137.57"
Set-TextValue "E26" "  -3.28%  "
Set-TextValue "E27" "  -12.19%  "
Set-TextValue "D28" "1.826"
Set-TextValue "E28" "  -9.59%  "
Set-TextValue "D29" "15.01"
Set-TextValue "E29" "  -9.09%  "
Set-TextValue "D30" "102.81"
Set-TextValue "E30" "  -5.62%  "
Set-TextValue "D31" "0.08203"
Set-TextValue "E31" "  -5.90%  "
Set-TextValue "D32" "3.687"
Set-TextValue "E32" "  -8.88%  "
Set-TextValue "D33" "3.456"
Set-TextValue "E33" "  -9.91%  "
Set-TextValue "D34" "0.04509"
Set-TextValue "E34" "  -3.67%  "
Set-TextValue "E35" "  -0.13%  "
Set-TextValue "D36" "2.659"
Set-TextValue "E36" "  -7.81%  "
Set-TextValue "D37" "0.9970"
Set-TextValue "E37" "  -9.57%  "
Set-TextValue "D38" "0.6073"
Set-TextValue "E38" "  -13.41%  "
Set-TextValue "D39" "2.696"
Set-TextValue "E39" "  -11.66%  "
Set-TextValue "B40" "VeChain"
Set-TextValue "C40" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D40" "0.01554"
Set-TextValue "E40" "  -4.76%  "
Set-TextValue "B41" "RenderToken"
Set-TextValue "C41" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D41" "1.944"
Set-TextValue "E41" "  -10.45%  "
Set-TextValue "E42" "  -0.20%  "
Set-TextValue "D43" "103.60"
Set-TextValue "E43" "  -1.27%  "
Set-TextValue "D44" "0.3818"
Set-TextValue "E44" "  -14.03%  "
Set-TextValue "D45" "0.7399"
Set-TextValue "E45" "  -13.76%  "
Set-TextValue "D46" "4.917"
Set-TextValue "E46" "  -14.14%  "
Set-TextValue "D47" "0.05470"
Set-TextValue "E47" "  -1.83%  "
Set-TextValue "D48" "0.1103"
Set-TextValue "E48" "  -5.55%  "
Set-TextValue "D49" "5.980"
Set-TextValue "E49" "  -14.20%  "
Set-TextValue "D50" "7.708"
Set-TextValue "E50" "  -10.05%  "
Set-TextValue "D51" "29.94"
Set-TextValue "E51" "  -10.06%  "
